$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("8059")

# --- Rename the sheet from "8059" to "4000" ---
$ws.Name = "4000"

# --- Update cell A2 (unit number) from 8059 to 4000 ---
# A2's number format is General, so a plain Range.Value assignment of a
# numeric-looking string ("4000") would be auto-coerced by Excel into a
# real number instead of staying text. Marking the cell as Text first keeps
# the written value a genuine string (A2 is also part of the merged range
# A2:D4 - switching NumberFormat, unlike Copy/PasteSpecial, does not disturb
# that merge).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "4000"

# --- Update FRA DEAD DATE from 01-06-2020 to 01-17-2020 (stays text) ---
# Same reasoning as A2: C5's number format ("M/D/YY;@") would otherwise
# cause "01-17-2020" to be auto-converted into a date serial number.
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "01-17-2020"

# --- Clear LUBE DUE (Y/N) value ---
$ws.Range("F5").Value = ""

# --- Update B24 from "AIR BRAKE" to "PTC TROUBLE" ---
$ws.Range("B24").Value = "PTC TROUBLE"

# --- Clear B25 (Alertor penalty, source still present) ---
$ws.Range("B25").Value = ""
